# Apply the "Add files via upload" changes to the Saldo/Export sheet:
#  - Add EMILIA (005535788) with Saldo 126143 (new highest balance)
#  - Correct LEONARDO's (004278212) balance from 12.76 to 10012.76 and move the
#    row to its new sorted position (the list is sorted by Saldo descending)
#  - Add MARIANA (005000460) with Saldo 3500

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the old LEONARDO (004278212) row that had the wrong (too small) balance.
#    This row currently sits far down the sheet (row 133) because 12.76 was a tiny value.
$ws.Rows.Item(133).Delete()

# 2) Insert MARIANA right before VALERIA (currently row 6) to keep the Saldo-descending order.
#    The leading apostrophe forces the account number to stay text (keeping the
#    leading zeros) while leaving the cell's number format as "General", just
#    like every other account-number cell already in the sheet.
$ws.Rows.Item(6).Insert()
$ws.Cells.Item(6,1).Value = "'005000460"
$ws.Cells.Item(6,2).Value = "MARIANA"
$ws.Cells.Item(6,3).Value = 3500

# 3) Insert the corrected LEONARDO row right before GUSTAVO (currently row 5).
$ws.Rows.Item(5).Insert()
$ws.Cells.Item(5,1).Value = "'004278212"
$ws.Cells.Item(5,2).Value = "LEONARDO"
$ws.Cells.Item(5,3).Value = 10012.76

# 4) Insert EMILIA right before ANA (currently row 2) - she now has the highest balance.
$ws.Rows.Item(2).Insert()
$ws.Cells.Item(2,1).Value = "'005535788"
$ws.Cells.Item(2,2).Value = "EMILIA"
$ws.Cells.Item(2,3).Value = 126143
